# Add WIN / TOP2 / TOP4 / RELEGATION columns between "Team" and "ExpPoints",
# pushing the existing ExpPoints column (and its values) from C to G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new, blank columns in front of the old "ExpPoints" column (C).
# This shifts the existing column C ("ExpPoints" + its values) to column G
# and carries the header style (bold/centered/bordered) along for the ride.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
